# Edit script for Result/VAR/Service/CHN.xlsx
# Applies corrected VAR/ARIMA fitted & predicted values across the four sheets.

$wb = $excel.ActiveWorkbook

$sheet1Data = @(
    "1981 22.46851474797015"
    "1982 23.34849440928966"
    "1983 23.81055562466604"
    "1984 24.84337344162781"
    "1985 27.09276523344505"
    "1986 30.42387063804104"
    "1987 30.97733424408571"
    "1988 31.39937535811302"
    "1989 32.19059041818007"
    "1990 33.56311052626297"
    "1991 33.29520472703416"
    "1992 34.67706599069223"
    "1993 35.54014274746507"
    "1994 34.86098136441871"
    "1995 35.1162105869378"
    "1996 34.96672565233457"
    "1997 35.47067191524393"
    "1998 36.7105286365936"
    "1999 38.17032523147132"
    "2000 39.29715163470701"
    "2001 40.27735997383937"
    "2002 41.12085486628393"
    "2003 41.79039250995109"
    "2004 41.68602196835528"
    "2005 41.36819002951953"
    "2006 41.95588640891205"
    "2007 42.68086585876535"
    "2008 43.71531677979664"
    "2009 43.97639471623226"
    "2010 44.07658620018655"
    "2011 43.83672058078202"
    "2012 45.47205014049865"
    "2013 46.92258918646102"
    "2014 48.6858898809519"
    "2015 50.46695266940336"
    "2016 52.59734367582512"
)
$sheet2Data = @(
    "49.32905490444891"
    "35.65399151601694"
    "-4.557407928155335"
    "-113.933392477008"
    "-405.8736453170176"
)
$sheet3Data = @(
    "1981 22.39727118482082"
    "1982 23.37063873846506"
    "1983 23.95391470777806"
    "1984 25.03882544703982"
    "1985 27.21138146027503"
    "1986 30.3491206092547"
    "1987 30.93927074559093"
    "1988 31.38926876876191"
    "1989 32.16460271319004"
    "1990 33.44712608422894"
    "1991 33.24326997548207"
    "1992 34.5073761242959"
    "1993 35.36631155174096"
    "1994 34.89566780412389"
    "1995 35.1526513954723"
    "1996 35.04932270534905"
    "1997 35.58545097780811"
    "1998 36.69720534650173"
    "1999 38.04504181654253"
    "2000 39.15811523831817"
    "2001 40.15506387977825"
    "2002 41.00775790188487"
    "2003 41.76132817779339"
    "2004 41.86336015031684"
    "2005 41.69781564042601"
    "2006 42.1797467789687"
    "2007 42.74719797823397"
    "2008 43.57215768122619"
    "2009 43.79980666805091"
    "2010 44.28402664963917"
    "2011 44.11420969486745"
    "2012 45.46164806268363"
    "2013 46.64096236423742"
    "2014 48.12840220237163"
    "2015 50.05394964943"
    "2016 52.66335510327869"
    "2017 53.15474797432245"
    "2018 53.68357589708526"
    "2019 54.40402564003569"
    "2020 54.07731666945551"
    "2021 53.43138774488179"
)
$sheet4Data = @(
    "52.22146493571091"
    "50.49310370692035"
    "48.4377000458832"
    "46.38145413436953"
    "44.60852370192259"
)

# --- Sheet 1: y_fitted_on_begin_2016 (dimension A1:D36 -> A1:D37) ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $parts = $sheet1Data[$i].Split(" ")
    $year = [double]$parts[0]
    $val = [double]$parts[1]
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $year
    $ws1.Cells.Item($row, 2).Value = $val
}
# New row 37 needs C/D "empty" cells matching the style of the row above it.
$ws1.Range("C36:D36").Copy($ws1.Range("C37:D37"))

# --- Sheet 2: y_pred_on_2017_2021 (dimension stays A1:D6, only column B changes) ---
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")
for ($i = 0; $i -lt $sheet2Data.Length; $i++) {
    $val = [double]$sheet2Data[$i]
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $val
}

# --- Sheet 3: y_fitted_on_begin_2021 (dimension A1:D43 -> A1:D42) ---
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
# Remove the last data row (row 43) first, since the new range has one fewer row.
$ws3.Rows.Item(43).Delete()
for ($i = 0; $i -lt $sheet3Data.Length; $i++) {
    $parts = $sheet3Data[$i].Split(" ")
    $year = [double]$parts[0]
    $val = [double]$parts[1]
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $year
    $ws3.Cells.Item($row, 2).Value = $val
}

# --- Sheet 4: y_pred_on_2022_2026 (dimension stays A1:D6, only column B changes) ---
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")
for ($i = 0; $i -lt $sheet4Data.Length; $i++) {
    $val = [double]$sheet4Data[$i]
    $row = $i + 2
    $ws4.Cells.Item($row, 2).Value = $val
}
